$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$ws.Range("B2").Value = 904380
$ws.Range("B3").Value = 823455
$ws.Range("B4").Value = 1345000
$ws.Range("B5").Value = 345700
$ws.Range("B6").Value = 465000

$ws.Range("G6").Select()
